$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 (01-01-2021) with revised figures
$ws.Range("B74").Value = 11966
$ws.Range("C74").Value = 8831
$ws.Range("D74").Value = 6235
$ws.Range("E74").Value = 1862
$ws.Range("F74").Value = 734
$ws.Range("G74").Value = 3788
$ws.Range("H74").Value = 196
$ws.Range("I74").Value = 3592
$ws.Range("J74").Value = 3075
$ws.Range("K74").Value = 518
$ws.Range("L74").Value = -3018
$ws.Range("M74").Value = 2366
$ws.Range("N74").Value = 1187
$ws.Range("O74").Value = 1092
$ws.Range("P74").Value = 87
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = 0

# Add new row 75 (01-04-2021)
# Force the cell to be treated as text so Excel doesn't auto-convert the
# dd-mm-yyyy-looking string into a date serial number, then clear the
# number-format override so the cell is left with no explicit style,
# matching the other "Serie" label cells in the column.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = 7286
$ws.Range("C75").Value = 4884
$ws.Range("D75").Value = 2970
$ws.Range("E75").Value = 2379
$ws.Range("F75").Value = -466
$ws.Range("G75").Value = 6358
$ws.Range("H75").Value = 811
$ws.Range("I75").Value = 5547
$ws.Range("J75").Value = 6426
$ws.Range("K75").Value = -880
$ws.Range("L75").Value = -2806
$ws.Range("M75").Value = -1150
$ws.Range("N75").Value = 326
$ws.Range("O75").Value = -1458
$ws.Range("P75").Value = -18
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = 0
